# UPDATE sections and lessons in course list XLSX file

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename "Sheet2" -> "Author" and populate it with the Author table
# ---------------------------------------------------------------------
$wsAuthor = $wb.Worksheets.Item("Sheet2")
$wsAuthor.Name = "Author"

$wsAuthor.Range("A1").Value = "AuthorID"
$wsAuthor.Range("B1").Value = "Author Name"
$wsAuthor.Range("A2").Value = 1
$wsAuthor.Range("B2").Value = "Bruce Myron"
$wsAuthor.Range("A3").Value = 1

# Column B needs to be wide enough to fit "Author Name" / "Bruce Myron"
$wsAuthor.Columns.Item(2).ColumnWidth = 11.92

# Matches the selection left behind in the saved file
[void]$wsAuthor.Range("B3").Select()

# ---------------------------------------------------------------------
# 2) Update the Section / Lesson / Date numbers on the Course sheet
# ---------------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("Course")

$wsCourse.Range("H2").Value = 5
$wsCourse.Range("I2").Value = 21
$wsCourse.Range("J2").Value = 45058

$wsCourse.Range("H3").Value = 4
$wsCourse.Range("I3").Value = 42
$wsCourse.Range("J3").Value = 45058

$wsCourse.Range("H4").Value = 2
$wsCourse.Range("I4").Value = 9

# Matches the selection left behind in the saved file
[void]$wsCourse.Range("C7").Select()
